$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark chapter 11 as done by moving the "START" marker to chapter 12
$ws.Range("B3").Value = "Ch 12 - START"
